# Add 8 more test cases for espn FantasyPage
$wb = $excel.ActiveWorkbook

# Rename the "Ebay" sheet to "URL" and populate it with the expected URLs
$urlSheet = $wb.Worksheets.Item("Ebay")
$urlSheet.Name = "URL"

$values = @(
    "expURL",
    "https://www.espn.com/fantasy/",
    "https://fantasy.espn.com/basketball/welcome?addata=fantasy_home_nav_fba2022",
    "https://fantasy.espn.com/hockey/welcome?addata=fhl_2022_fantasy_home_nav",
    "https://www.espn.com/fantasy/baseball/",
    "https://www.espn.com/fantasy/football/",
    "https://fantasy.espn.com/free-prize-games",
    "https://www.espn.com/fantasy/basketball/",
    "https://www.espn.com/fantasy/hockey/",
    "https://fantasy.espn.com/streak/en/",
    "http://www.espn.com/espn/apps/fantasy"
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $urlSheet.Cells.Item($row, 1).Value = $values[$i]
}

$urlSheet.Columns.Item(1).ColumnWidth = 70.83
[void]$urlSheet.Range("A3").Select()

# Adjust the NavBarMenu sheet: widen column A and change selection
$navSheet = $wb.Worksheets.Item("NavBarMenu")
$navSheet.Columns.Item(1).ColumnWidth = 39.33
[void]$navSheet.Range("B1:B1048576").Select()
$navSheet.PageSetup.Orientation = 1

# Make URL sheet the active sheet/tab
[void]$urlSheet.Activate()

$wb.Save()
